$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for Price/Volume columns so Excel does not
# auto-coerce values like "53.90" or "6.420" into numbers (which
# would silently drop the significant trailing/format digits).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.924.66"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.705.76"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.92"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4026"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4077"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.57%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.90"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.473"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08829"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.21"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.496"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.043"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001350"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.673.31"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.68"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07176"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.272"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.51"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.924.36"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.334"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.883"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.420"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +22.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.14"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.54"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "143.84"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.217"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.270"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +14.03%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.839.87"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.01%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08735"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.98%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.03205"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +8.99%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.363"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.036"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2869"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8517"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +7.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.86"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09457"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.68%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.473"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.64"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.729"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7485"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.239"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.392"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.001"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.32"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.08409"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.38%  "
